$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 6.654043666666666
$ws.Cells.Item(2, 8).Value = 19.962131
$ws.Cells.Item(2, 9).Value = 0.3091924566209486
$ws.Cells.Item(2, 10).Value = 0.3091924566209486
$ws.Cells.Item(2, 13).Value = 6.956267333333333
$ws.Cells.Item(2, 14).Value = 20.868802
$ws.Cells.Item(2, 15).Value = 0.383240417447883
$ws.Cells.Item(2, 16).Value = 0.383240417447883
$ws.Cells.Item(2, 17).Value = 46.28730659300688
$ws.Cells.Item(2, 18).Value = 416.585759337062
$ws.Cells.Item(2, 19).Value = 0.1184950461471488
$ws.Cells.Item(2, 20).Value = 0.1184950461471488
$ws.Cells.Item(3, 7).Value = 6.654043666666666
$ws.Cells.Item(3, 8).Value = 19.962131
$ws.Cells.Item(3, 9).Value = 0.3091924566209486
$ws.Cells.Item(3, 10).Value = 0.3091924566209486
$ws.Cells.Item(3, 15).Value = 0.3291757349456286
$ws.Cells.Item(3, 16).Value = 0.3291757349456286
$ws.Cells.Item(3, 17).Value = 39.75744068924755
$ws.Cells.Item(3, 18).Value = 357.816966203228
$ws.Cells.Item(3, 19).Value = 0.1017786541478451
$ws.Cells.Item(3, 20).Value = 0.1017786541478451
$ws.Cells.Item(4, 7).Value = 6.654043666666666
$ws.Cells.Item(4, 8).Value = 19.962131
$ws.Cells.Item(4, 9).Value = 0.3091924566209486
$ws.Cells.Item(4, 10).Value = 0.3091924566209486
$ws.Cells.Item(4, 13).Value = 3.345755333333333
$ws.Cells.Item(4, 14).Value = 10.037266
$ws.Cells.Item(4, 15).Value = 0.1843271123984713
$ws.Cells.Item(4, 16).Value = 0.1843271123984713
$ws.Cells.Item(4, 17).Value = 22.26280208598289
$ws.Cells.Item(4, 18).Value = 200.365218773846
$ws.Cells.Item(4, 19).Value = 0.05699255270432905
$ws.Cells.Item(4, 20).Value = 0.05699255270432905
$ws.Cells.Item(5, 7).Value = 6.654043666666666
$ws.Cells.Item(5, 8).Value = 19.962131
$ws.Cells.Item(5, 9).Value = 0.3091924566209486
$ws.Cells.Item(5, 10).Value = 0.3091924566209486
$ws.Cells.Item(5, 13).Value = 1.874232
$ws.Cells.Item(5, 14).Value = 5.622696
$ws.Cells.Item(5, 15).Value = 0.1032567352080173
$ws.Cells.Item(5, 16).Value = 0.1032567352080173
$ws.Cells.Item(5, 17).Value = 12.471221569464
$ws.Cells.Item(5, 18).Value = 112.240994125176
$ws.Cells.Item(5, 19).Value = 0.03192620362162567
$ws.Cells.Item(5, 20).Value = 0.03192620362162566
$ws.Cells.Item(6, 9).Value = 0.09233579784218476
$ws.Cells.Item(6, 10).Value = 0.09233579784218476
$ws.Cells.Item(6, 13).Value = 6.956267333333333
$ws.Cells.Item(6, 14).Value = 20.868802
$ws.Cells.Item(6, 15).Value = 0.383240417447883
$ws.Cells.Item(6, 16).Value = 0.383240417447883
$ws.Cells.Item(6, 17).Value = 13.82302605613289
$ws.Cells.Item(6, 18).Value = 124.407234505196
$ws.Cells.Item(6, 19).Value = 0.03538680971042222
$ws.Cells.Item(6, 20).Value = 0.03538680971042222
$ws.Cells.Item(7, 9).Value = 0.09233579784218476
$ws.Cells.Item(7, 10).Value = 0.09233579784218476
$ws.Cells.Item(7, 15).Value = 0.3291757349456286
$ws.Cells.Item(7, 16).Value = 0.3291757349456286
$ws.Cells.Item(7, 19).Value = 0.03039470411649215
$ws.Cells.Item(7, 20).Value = 0.03039470411649215
$ws.Cells.Item(8, 9).Value = 0.09233579784218476
$ws.Cells.Item(8, 10).Value = 0.09233579784218476
$ws.Cells.Item(8, 13).Value = 3.345755333333333
$ws.Cells.Item(8, 14).Value = 10.037266
$ws.Cells.Item(8, 15).Value = 0.1843271123984713
$ws.Cells.Item(8, 16).Value = 0.1843271123984713
$ws.Cells.Item(8, 17).Value = 6.648459717540889
$ws.Cells.Item(8, 18).Value = 59.83613745786801
$ws.Cells.Item(8, 19).Value = 0.01701999098725892
$ws.Cells.Item(8, 20).Value = 0.01701999098725892
$ws.Cells.Item(9, 9).Value = 0.09233579784218476
$ws.Cells.Item(9, 10).Value = 0.09233579784218476
$ws.Cells.Item(9, 13).Value = 1.874232
$ws.Cells.Item(9, 14).Value = 5.622696
$ws.Cells.Item(9, 15).Value = 0.1032567352080173
$ws.Cells.Item(9, 16).Value = 0.1032567352080173
$ws.Cells.Item(9, 17).Value = 3.724347632112
$ws.Cells.Item(9, 18).Value = 33.519128689008
$ws.Cells.Item(9, 19).Value = 0.009534293028011488
$ws.Cells.Item(9, 20).Value = 0.009534293028011486
$ws.Cells.Item(10, 7).Value = 0.9593116666666667
$ws.Cells.Item(10, 8).Value = 2.877935
$ws.Cells.Item(10, 9).Value = 0.04457619242381535
$ws.Cells.Item(10, 10).Value = 0.04457619242381536
$ws.Cells.Item(10, 13).Value = 6.956267333333333
$ws.Cells.Item(10, 14).Value = 20.868802
$ws.Cells.Item(10, 15).Value = 0.383240417447883
$ws.Cells.Item(10, 16).Value = 0.383240417447883
$ws.Cells.Item(10, 17).Value = 6.673228409318889
$ws.Cells.Item(10, 18).Value = 60.05905568386999
$ws.Cells.Item(10, 19).Value = 0.01708339859274015
$ws.Cells.Item(10, 20).Value = 0.01708339859274016
$ws.Cells.Item(11, 7).Value = 0.9593116666666667
$ws.Cells.Item(11, 8).Value = 2.877935
$ws.Cells.Item(11, 9).Value = 0.04457619242381535
$ws.Cells.Item(11, 10).Value = 0.04457619242381536
$ws.Cells.Item(11, 15).Value = 0.3291757349456286
$ws.Cells.Item(11, 16).Value = 0.3291757349456286
$ws.Cells.Item(11, 17).Value = 5.731819416975555
$ws.Cells.Item(11, 18).Value = 51.58637475278
$ws.Cells.Item(11, 19).Value = 0.01467340090218718
$ws.Cells.Item(11, 20).Value = 0.01467340090218718
$ws.Cells.Item(12, 7).Value = 0.9593116666666667
$ws.Cells.Item(12, 8).Value = 2.877935
$ws.Cells.Item(12, 9).Value = 0.04457619242381535
$ws.Cells.Item(12, 10).Value = 0.04457619242381536
$ws.Cells.Item(12, 13).Value = 3.345755333333333
$ws.Cells.Item(12, 14).Value = 10.037266
$ws.Cells.Item(12, 15).Value = 0.1843271123984713
$ws.Cells.Item(12, 16).Value = 0.1843271123984713
$ws.Cells.Item(12, 17).Value = 3.209622125078889
$ws.Cells.Item(12, 18).Value = 28.88659912571
$ws.Cells.Item(12, 19).Value = 0.008216600831200497
$ws.Cells.Item(12, 20).Value = 0.008216600831200498
$ws.Cells.Item(13, 7).Value = 0.9593116666666667
$ws.Cells.Item(13, 8).Value = 2.877935
$ws.Cells.Item(13, 9).Value = 0.04457619242381535
$ws.Cells.Item(13, 10).Value = 0.04457619242381536
$ws.Cells.Item(13, 13).Value = 1.874232
$ws.Cells.Item(13, 14).Value = 5.622696
$ws.Cells.Item(13, 15).Value = 0.1032567352080173
$ws.Cells.Item(13, 16).Value = 0.1032567352080173
$ws.Cells.Item(13, 17).Value = 1.79797262364
$ws.Cells.Item(13, 18).Value = 16.18175361276
$ws.Cells.Item(13, 19).Value = 0.00460279209768753
$ws.Cells.Item(13, 20).Value = 0.00460279209768753
$ws.Cells.Item(14, 7).Value = 11.92023
$ws.Cells.Item(14, 8).Value = 35.76069
$ws.Cells.Item(14, 9).Value = 0.5538955531130513
$ws.Cells.Item(14, 10).Value = 0.5538955531130514
$ws.Cells.Item(14, 13).Value = 6.956267333333333
$ws.Cells.Item(14, 14).Value = 20.868802
$ws.Cells.Item(14, 15).Value = 0.383240417447883
$ws.Cells.Item(14, 16).Value = 0.383240417447883
$ws.Cells.Item(14, 17).Value = 82.92030655481999
$ws.Cells.Item(14, 18).Value = 746.2827589933798
$ws.Cells.Item(14, 19).Value = 0.2122751629975718
$ws.Cells.Item(14, 20).Value = 0.2122751629975718
$ws.Cells.Item(15, 7).Value = 11.92023
$ws.Cells.Item(15, 8).Value = 35.76069
$ws.Cells.Item(15, 9).Value = 0.5538955531130513
$ws.Cells.Item(15, 10).Value = 0.5538955531130514
$ws.Cells.Item(15, 15).Value = 0.3291757349456286
$ws.Cells.Item(15, 16).Value = 0.3291757349456286
$ws.Cells.Item(15, 17).Value = 71.22253188707998
$ws.Cells.Item(15, 18).Value = 641.0027869837199
$ws.Cells.Item(15, 19).Value = 0.1823289757791041
$ws.Cells.Item(15, 20).Value = 0.1823289757791041
$ws.Cells.Item(16, 7).Value = 11.92023
$ws.Cells.Item(16, 8).Value = 35.76069
$ws.Cells.Item(16, 9).Value = 0.5538955531130513
$ws.Cells.Item(16, 10).Value = 0.5538955531130514
$ws.Cells.Item(16, 13).Value = 3.345755333333333
$ws.Cells.Item(16, 14).Value = 10.037266
$ws.Cells.Item(16, 15).Value = 0.1843271123984713
$ws.Cells.Item(16, 16).Value = 0.1843271123984713
$ws.Cells.Item(16, 17).Value = 39.88217309706
$ws.Cells.Item(16, 18).Value = 358.93955787354
$ws.Cells.Item(16, 19).Value = 0.1020979678756828
$ws.Cells.Item(16, 20).Value = 0.1020979678756829
$ws.Cells.Item(17, 7).Value = 11.92023
$ws.Cells.Item(17, 8).Value = 35.76069
$ws.Cells.Item(17, 9).Value = 0.5538955531130513
$ws.Cells.Item(17, 10).Value = 0.5538955531130514
$ws.Cells.Item(17, 13).Value = 1.874232
$ws.Cells.Item(17, 14).Value = 5.622696
$ws.Cells.Item(17, 15).Value = 0.1032567352080173
$ws.Cells.Item(17, 16).Value = 0.1032567352080173
$ws.Cells.Item(17, 17).Value = 22.34127651336
$ws.Cells.Item(17, 18).Value = 201.07148862024
$ws.Cells.Item(17, 19).Value = 0.05719344646069262
$ws.Cells.Item(17, 20).Value = 0.05719344646069263
